$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the revision-history table row for "EDUARDO DORADOR PINA"
#    (whole <w:tr> row, including the vMerge continuation cells)
# ------------------------------------------------------------------
$t = $d.Tables.Item(1)
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $row = $t.Rows.Item($i)
    if ($row.Cells.Item(1).Range.Text -like "*EDUARDO DORADOR PINA*") {
        $row.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2) Change " e progresso do mesmo" -> " e progresso " + "dele"
#    The target keeps " e progresso " in the existing run and adds a
#    brand-new run (identical formatting) that just contains "dele".
#    Turning TrackRevisions on for this one edit, then accepting only
#    the two revisions it produces, forces Word to keep the text split
#    across two separate <w:r> elements instead of silently
#    re-merging them into a single run while saving.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" e progresso do mesmo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$oldText = " e progresso do mesmo"
$splitAt = $oldText.IndexOf("do mesmo")
$subStart = $rng.Start + $splitAt
$subEnd = $rng.End
$sub = $d.Range($subStart, $subEnd)

$d.TrackRevisions = $true
$sub.Text = "dele"
$d.TrackRevisions = $false

for ($k = $d.Revisions.Count; $k -ge 1; $k--) {
    $d.Revisions.Item($k).Accept()
}
